# Actualizacion Datos Personales 4 nov
#
# The group "1AM" (previously taught by "Saucedo Rivalcoba Graciela", stored in the
# last data row of each parcial sheet) is now taught by "García Sánchez Magda Bexabe",
# and the "1AM" row is relocated so it sits right after the "Contreras Díaz" block
# (i.e. immediately before "1BM"), pushing the other "García Sánchez" rows
# (1BM, 1CM, 1DM, 1EM, 1FM) down by one row.

$wb = $excel.ActiveWorkbook

$cols = @("A","B","C","D","E","F","G","H","I","J","K")

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Remember the last row's (row 12 - group 1AM / Saucedo Rivalcoba Graciela) data
    # before it gets overwritten by the downward shift below.
    $lastRow = @{}
    foreach ($col in $cols) {
        $lastRow[$col] = $ws.Range($col + "12").Value()
    }

    # Shift rows 7-11 down into rows 8-12 (bottom-up so we don't clobber data
    # before it has been read).
    for ($r = 11; $r -ge 7; $r--) {
        foreach ($col in $cols) {
            $ws.Range($col + ($r + 1)).Value = $ws.Range($col + $r).Value()
        }
    }

    # Row 7 now becomes the former last row (group 1AM), re-assigned to teacher
    # "García Sánchez Magda Bexabe" instead of "Saucedo Rivalcoba Graciela".
    $ws.Range("A7").Value = "García Sánchez Magda Bexabe"
    $ws.Range("B7").Value = "1AM"
    foreach ($col in @("C","D","E","F","G","H","I","J","K")) {
        $ws.Range($col + "7").Value = $lastRow[$col]
    }
}
